$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Copy the formatting (cell styles) of the legend/comparison table (which
#    currently lives at K6:S10) down onto its new home at K13:R17, BEFORE we
#    touch any of the source cells. We copy in the exact groups that line up
#    with the merged-cell layout of the destination table.
# ---------------------------------------------------------------------------

# Row 13 (legend header band) <- old row 6 (L6:S6)
$ws.Range("L6:M6").Copy()
$ws.Range("L13").PasteSpecial(-4122)
$ws.Range("O6:P6").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("Q6:S6").Copy()
$ws.Range("P13").PasteSpecial(-4122)

# Row 14 (column headers) <- old row 7 (K7:S7), dropping the old M7 (OCR) column
$ws.Range("K7").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("L7").Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("N7").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("O7:P7").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("Q7:S7").Copy()
$ws.Range("P14").PasteSpecial(-4122)

# Row 15 (ATRT) <- old row 8 (K8:S8), dropping old M8
$ws.Range("K8").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("L8").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("N8").Copy()
$ws.Range("M15").PasteSpecial(-4122)
$ws.Range("O8:P8").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("Q8:S8").Copy()
$ws.Range("P15").PasteSpecial(-4122)

# Row 16 (eggPlant) <- old row 9 (K9:S9), dropping old M9
$ws.Range("K9").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("L9").Copy()
$ws.Range("L16").PasteSpecial(-4122)
$ws.Range("N9").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("O9:P9").Copy()
$ws.Range("N16").PasteSpecial(-4122)
$ws.Range("Q9:S9").Copy()
$ws.Range("P16").PasteSpecial(-4122)

# Row 17 (Sikuli) <- old row 10 (K10:S10), dropping old M10
$ws.Range("K10").Copy()
$ws.Range("K17").PasteSpecial(-4122)
$ws.Range("L10").Copy()
$ws.Range("L17").PasteSpecial(-4122)
$ws.Range("N10").Copy()
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("O10:P10").Copy()
$ws.Range("N17").PasteSpecial(-4122)
$ws.Range("Q10:S10").Copy()
$ws.Range("P17").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Clear out the old K6:S10 block now that its formatting has been copied.
# ---------------------------------------------------------------------------
$ws.Range("L6:S6").ClearContents()
$ws.Range("K7:S10").ClearContents()

# ---------------------------------------------------------------------------
# 3. Un-merge the old legend merges and merge the new ones.
# ---------------------------------------------------------------------------
$ws.Range("L6:N6").UnMerge()
$ws.Range("O6:P6").UnMerge()
$ws.Range("Q6:S6").UnMerge()

$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()
$ws.Range("P13:R13").Merge()

# ---------------------------------------------------------------------------
# 4. Fill in the values of the new comparison table (K13:R17).
# ---------------------------------------------------------------------------
$ws.Range("L13").Value = "SUT Interaction and Performance"
$ws.Range("N13").Value = "Collaboration"
$ws.Range("P13").Value = "Other"

$ws.Range("K14").Value = "Tool"
$ws.Range("L14").Value = "Image Capture and Scan"
$ws.Range("M14").Value = "Time To Execute Common Scenario"
$ws.Range("N14").Value = "IDE Source File Creation And Management"
$ws.Range("O14").Value = "Source Revision Control"
$ws.Range("P14").Value = "Test Execution Reporting Capabilities"
$ws.Range("Q14").Value = "Linking Requirements To Test Steps"
$ws.Range("R14").Value = "License"

$ws.Range("K15").Value = "ATRT `nVersion 5.6.8"
$ws.Range("L15").Value = 8
$ws.Range("M15").Value = "Success: 2 min 45 sec`nFailure: NOT TESTED due to management decision"
$ws.Range("N15").Value = 4
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 8
$ws.Range("Q15").Value = "Yes but not tested"
$ws.Range("R15").Value = "Commercial.  See Appendix A"

$ws.Range("K16").Value = "eggPlant`nVersion 14.01"
$ws.Range("L16").Value = 8
$ws.Range("M16").Value = "Success: 55 secs`nFailure: 40 secs"
$ws.Range("N16").Value = 7
$ws.Range("O16").Value = 8
$ws.Range("P16").Value = 7
$ws.Range("Q16").Value = "No"
$ws.Range("R16").Value = "Commercial.  See Appendix A"

$ws.Range("K17").Value = "Sikuli`nVersion 1.0.1"
$ws.Range("L17").Value = 8
$ws.Range("M17").Value = "Success: 2 min 20 sec`nFailure: 1 min 20 sec"
$ws.Range("N17").Value = 7
$ws.Range("O17").Value = 8
$ws.Range("P17").Value = 7
$ws.Range("Q17").Value = "No"
$ws.Range("R17").Value = "Open Source `nMIT License"

# ---------------------------------------------------------------------------
# 5. Update the first (tool-comparison) table.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 47.25
$ws.Rows.Item(8).RowHeight = 75

$ws.Range("I8").Value = "6hrs 30mins for success scenario only.  Failure scenario not finished"
$ws.Range("H9").Value = 8
$ws.Range("I9").Value = "4 hours for success and failure scenarios"
$ws.Range("I10").Value = "5hrs 45mins for success and failure scenarios"

# ---------------------------------------------------------------------------
# 6. Column widths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 17.95
$ws.Columns.Item(12).ColumnWidth = 12.6
$ws.Columns.Item(13).ColumnWidth = 22.6
$ws.Columns.Item(14).ColumnWidth = 17.45
$ws.Columns.Item(15).ColumnWidth = 16.6
$ws.Columns.Item(16).ColumnWidth = 12.1
$ws.Columns.Item(17).ColumnWidth = 15.1
$ws.Columns.Item(18).ColumnWidth = 12.75

# ---------------------------------------------------------------------------
# 7. Sheet view (scroll position / selection).
# ---------------------------------------------------------------------------
$ws.Range("M10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 9
